$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.369083150264215
$ws.Range("D2").Value = 4.181045905317358
$ws.Range("E2").Value = 16.52620225244815
$ws.Range("F2").Value = 21.82287381936986
$ws.Range("G2").Value = 3.608441093530877
$ws.Range("I2").Value = 18.6212245218672
$ws.Range("K2").Value = 11.83129128825535
$ws.Range("O2").Value = 19.40173363499445
$ws.Range("B3").Value = 6.222685135686072
$ws.Range("D3").Value = 4.139874056522111
$ws.Range("E3").Value = 15.5865753296442
$ws.Range("F3").Value = 21.81140239631509
$ws.Range("G3").Value = 3.61066802048957
$ws.Range("I3").Value = 18.74925458363776
$ws.Range("K3").Value = 11.24110303337065
$ws.Range("O3").Value = 19.4607630822676
$ws.Range("B4").Value = 6.131641048875076
$ws.Range("D4").Value = 4.114164890604886
$ws.Range("E4").Value = 14.98470727304621
$ws.Range("F4").Value = 21.81276049557438
$ws.Range("G4").Value = 3.612105768788702
$ws.Range("I4").Value = 18.83181888369029
$ws.Range("K4").Value = 10.86039897957681
$ws.Range("O4").Value = 19.5034544111015
$ws.Range("B5").Value = 6.094299572387109
$ws.Range("D5").Value = 4.103586057230975
$ws.Range("E5").Value = 14.73344042902481
$ws.Range("F5").Value = 21.81542157554045
$ws.Range("G5").Value = 3.612709428818296
$ws.Range("I5").Value = 18.866461644705
$ws.Range("K5").Value = 10.7007553308877
$ws.Range("O5").Value = 19.52246301368287
$ws.Range("B6").Value = 6.088086046287505
$ws.Range("D6").Value = 4.101823453669972
$ws.Range("E6").Value = 14.69136424027861
$ws.Range("F6").Value = 21.81599055238518
$ws.Range("G6").Value = 3.612810740937765
$ws.Range("I6").Value = 18.87227436308498
$ws.Range("K6").Value = 10.67397806386422
$ws.Range("O6").Value = 19.52571641978081
$ws.Range("B7").Value = 6.131138353582328
$ws.Range("D7").Value = 4.114022626547362
$ws.Range("E7").Value = 14.98134250279656
$ws.Range("F7").Value = 21.81278785841819
$ws.Range("G7").Value = 3.612113837943497
$ws.Range("I7").Value = 18.83228204668042
$ws.Range("K7").Value = 10.85826405189587
$ws.Range("O7").Value = 19.50370425586435
$ws.Range("B8").Value = 6.318874273507927
$ws.Range("D8").Value = 4.166942478099771
$ws.Range("E8").Value = 16.20753948384126
$ws.Range("F8").Value = 21.81717238336799
$ws.Range("G8").Value = 3.609194362396368
$ws.Range("I8").Value = 18.66455061025835
$ws.Range("K8").Value = 11.63165478656741
$ws.Range("O8").Value = 19.42074328290598
$ws.Range("B9").Value = 6.675772529890628
$ws.Range("D9").Value = 4.267058598531326
$ws.Range("E9").Value = 18.501013184956
$ws.Range("F9").Value = 21.8925522302117
$ws.Range("G9").Value = 3.604025097945001
$ws.Range("I9").Value = 18.36685798813219
$ws.Range("K9").Value = 12.99918838720052
$ws.Range("O9").Value = 19.30961761499689
$ws.Range("B10").Value = 6.928516286912388
$ws.Range("D10").Value = 4.33805449058059
$ws.Range("E10").Value = 20.14724419895321
$ws.Range("F10").Value = 21.98865524229661
$ws.Range("G10").Value = 3.600562120330933
$ws.Range("I10").Value = 18.16698419320571
$ws.Range("K10").Value = 13.90918758900551
$ws.Range("O10").Value = 19.25992583507372
$ws.Range("B11").Value = 7.040943287505756
$ws.Range("D11").Value = 4.369729366127005
$ws.Range("E11").Value = 20.85394232886489
$ws.Range("F11").Value = 22.04116942729239
$ws.Range("G11").Value = 3.59905859290232
$ws.Range("I11").Value = 18.08010579589082
$ws.Range("K11").Value = 14.30203965011906
$ws.Range("O11").Value = 19.24435615835242
$ws.Range("B12").Value = 7.083112683071707
$ws.Range("D12").Value = 4.381628953372846
$ws.Range("E12").Value = 21.1155218924709
$ws.Range("F12").Value = 22.06231220509725
$ws.Range("G12").Value = 3.598499506276694
$ws.Range("I12").Value = 18.04778578141489
$ws.Range("K12").Value = 14.44773232748647
$ws.Range("O12").Value = 19.23947896156193
$ws.Range("B13").Value = 7.074049302337805
$ws.Range("D13").Value = 4.379070477407445
$ws.Range("E13").Value = 21.05945377731882
$ws.Range("F13").Value = 22.05770298793009
$ws.Range("G13").Value = 3.598619459894441
$ws.Range("I13").Value = 18.05472076929073
$ws.Range("K13").Value = 14.41649191560178
$ws.Range("O13").Value = 19.24048394216189
$ws.Range("B14").Value = 7.044420898180396
$ws.Range("D14").Value = 4.370710287948438
$ws.Range("E14").Value = 20.87558332781265
$ws.Range("F14").Value = 22.04288372329088
$ws.Range("G14").Value = 3.599012391082378
$ws.Range("I14").Value = 18.07743522146741
$ws.Range("K14").Value = 14.31408761775771
$ws.Range("O14").Value = 19.24393445454349
$ws.Range("B15").Value = 7.026218921016674
$ws.Range("D15").Value = 4.365576894339926
$ws.Range("E15").Value = 20.76217286599086
$ws.Range("F15").Value = 22.03396988590145
$ws.Range("G15").Value = 3.599254408077416
$ws.Range("I15").Value = 18.09142379457907
$ws.Range("K15").Value = 14.25096104528357
$ws.Range("O15").Value = 19.24618085387827
$ws.Range("B16").Value = 6.921114345949092
$ws.Range("D16").Value = 4.335971471806712
$ws.Range("E16").Value = 20.10021393021887
$ws.Range("F16").Value = 21.98539968406904
$ws.Range("G16").Value = 3.600661819192049
$ws.Range("I16").Value = 18.17274306512958
$ws.Range("K16").Value = 13.88308569831544
$ws.Range("O16").Value = 19.26108553693826
$ws.Range("B17").Value = 6.85595570148731
$ws.Range("D17").Value = 4.317646265189416
$ws.Range("E17").Value = 19.68334424546295
$ws.Range("F17").Value = 21.95785140120844
$ws.Range("G17").Value = 3.601543568301298
$ws.Range("I17").Value = 18.22366391273087
$ws.Range("K17").Value = 13.65197067426188
$ws.Range("O17").Value = 19.27203623416185
$ws.Range("B18").Value = 6.818240480913349
$ws.Range("D18").Value = 4.307047918156818
$ws.Range("E18").Value = 19.43959958096881
$ws.Range("F18").Value = 21.9428349598142
$ws.Range("G18").Value = 3.60205748852097
$ws.Range("I18").Value = 18.25333315248541
$ws.Range("K18").Value = 13.5170539316951
$ws.Range("O18").Value = 19.27899661250244
$ws.Range("B19").Value = 6.805431136269722
$ws.Range("D19").Value = 4.303449682826948
$ws.Range("E19").Value = 19.35638855338162
$ws.Range("F19").Value = 21.937893172327
$ws.Range("G19").Value = 3.602232656134023
$ws.Range("I19").Value = 18.26344415395171
$ws.Range("K19").Value = 13.47103366628092
$ws.Range("O19").Value = 19.28146672041185
$ws.Range("B20").Value = 6.862916834909637
$ws.Range("D20").Value = 4.319603076647803
$ws.Range("E20").Value = 19.72813160204963
$ws.Range("F20").Value = 21.96069825688308
$ws.Range("G20").Value = 3.60144900523858
$ws.Range("I20").Value = 18.21820389824815
$ws.Range("K20").Value = 13.67677904158615
$ws.Range("O20").Value = 19.27080196383158
$ws.Range("B21").Value = 7.053134738360152
$ws.Range("D21").Value = 4.373168500584478
$ws.Range("E21").Value = 20.9297539654134
$ws.Range("F21").Value = 22.04720246801111
$ws.Range("G21").Value = 3.59889669949615
$ws.Range("I21").Value = 18.0707477454758
$ws.Range("K21").Value = 14.34424986031889
$ws.Range("O21").Value = 19.24289325583557
$ws.Range("B22").Value = 7.175081200908864
$ws.Range("D22").Value = 4.407620370835025
$ws.Range("E22").Value = 21.679941937827
$ws.Range("F22").Value = 22.11105809836445
$ws.Range("G22").Value = 3.5972884349617
$ws.Range("I22").Value = 17.97774993069459
$ws.Range("K22").Value = 14.76256402770261
$ws.Range("O22").Value = 19.23059380309178
$ws.Range("B23").Value = 7.110224892778916
$ws.Range("D23").Value = 4.389285471412766
$ws.Range("E23").Value = 21.28275684322951
$ws.Range("F23").Value = 22.07631065132975
$ws.Range("G23").Value = 3.598141342020423
$ws.Range("I23").Value = 18.0270768564363
$ws.Range("K23").Value = 14.54095117753526
$ws.Range("O23").Value = 19.23661256914736
$ws.Range("B24").Value = 6.859770496763884
$ws.Range("D24").Value = 4.31871859777553
$ws.Range("E24").Value = 19.70789595202149
$ws.Range("F24").Value = 21.95940863401739
$ws.Range("G24").Value = 3.601491735417175
$ws.Range("I24").Value = 18.22067114273992
$ws.Range("K24").Value = 13.6655695477242
$ws.Range("O24").Value = 19.27135790743702
$ws.Range("B25").Value = 6.580706845584647
$ws.Range("D25").Value = 4.240400674237723
$ws.Range("E25").Value = 17.8571438889858
$ws.Range("F25").Value = 21.86500014378401
$ws.Range("G25").Value = 3.605364427373635
$ws.Range("I25").Value = 18.44406896116215
$ws.Range("K25").Value = 12.64564317716832
$ws.Range("O25").Value = 19.33410408350478
